# Apply the "Add data for 2022-10-19" update to the carjacking-by-neighborhood
# workbook: rename the reporting sheet / header label from "October 10" to
# "October 11", and update the underlying counts for the affected
# neighborhood/day-of-month columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab and update the header label in B1 to reflect the
# new "through" date.
$ws.Name = "Through 2022-10-11"
$ws.Range("B1").Value = "October 2022 (through October 11)"

# Cell updates: Range address -> new value
$updates = @{
    "B2"   = 2
    "L2"   = 8
    "BT2"  = 3
    "BJ3"  = 3
    "AP4"  = 1
    "L6"   = 6
    "AP6"  = 1
    "L9"   = 1
    "B13"  = 3
    "V18"  = 4
    "V19"  = 2
    "B23"  = 1
    "L32"  = 5
    "AZ38" = 1
    "L48"  = 1
    "V52"  = 1
    "V78"  = 1
    "V80"  = 2
    "AP80" = 3
    "B93"  = 1
    "V95"  = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
